# "update to version 0.9"
#
# The sheet tracked two separate death-rate rows ("WildTypeDeath" and
# "MutantDeath"); this revision collapses them into a single "Death" row.
#
# Concretely: row 9 ("MutantDeath") is removed entirely (shifting every
# row below it up by one), and the label that used to read "WildTypeDeath"
# (now sitting in row 8) is renamed to simply "Death".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole row 9 - this shifts rows 10..40 up to 9..39 and updates
# the sheet's used range/dimension automatically.
$ws.Rows(9).Delete()

# Row 8 (formerly "WildTypeDeath") becomes just "Death".
$ws.Range("A8").Value = "Death"

# Leave the selection where the author left it when they saved.
$ws.Range("A8").Select()
